# Repull data: update column F (dSF) values for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -1
    6  = 2
    8  = 0
    10 = 0
    12 = 1
    13 = -1
    15 = 1
    18 = -1
    25 = -2
    43 = -1
    44 = 2
    45 = -1
    47 = 2
    48 = 2
    53 = 1
    55 = 2
    58 = 0
    63 = 1
    68 = 6
    73 = -3
    76 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
